$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in column D stay as text by forcing
# the cell's number format to Text before assigning the string value.
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.355.58'
$ws.Range('E2').Value = '  +1.74%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.865.57'
$ws.Range('E4').Value = '  +1.36%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '314.24'
$ws.Range('E5').Value = '  +1.54%  '
$ws.Range('E6').Value = '  +1.27%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4809'
$ws.Range('E7').Value = '  +2.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3734'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07424'
$ws.Range('E9').Value = '  +3.76%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9384'
$ws.Range('E10').Value = '  +2.02%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.78'
$ws.Range('E11').Value = '  +6.46%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07892'
$ws.Range('E12').Value = '  +3.73%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.863.00'
$ws.Range('E13').Value = '  +0.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.441'
$ws.Range('E14').Value = '  +3.02%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.558'
$ws.Range('E15').Value = '  +2.57%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.42'
$ws.Range('E16').Value = '  +2.95%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.023'
$ws.Range('E17').Value = '  +1.33%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008801'
$ws.Range('E18').Value = '  +1.99%  '
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '27.384.48'
$ws.Range('E21').Value = '  +1.77%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.73'
$ws.Range('E23').Value = '  +1.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.962'
$ws.Range('E24').Value = '  +1.81%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.32'
$ws.Range('E25').Value = '  +1.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.56'
$ws.Range('E26').Value = '  +2.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.012'
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '116.14'
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.002'
$ws.Range('E29').Value = '  +3.00%  '
$ws.Range('E30').Value = '  +1.30%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.353'
$ws.Range('E31').Value = '  +3.79%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.198'
$ws.Range('E32').Value = '  +2.47%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.577'
$ws.Range('E33').Value = '  +2.29%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7451'
$ws.Range('E34').Value = '  +0.41%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.697'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02055'
$ws.Range('E36').Value = '  +5.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.126'
$ws.Range('E37').Value = '  +3.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05304'
$ws.Range('E38').Value = '  +1.38%  '
$ws.Range('B39').Value = 'TheSandbox'
$ws.Range('C39').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.5398'
$ws.Range('E39').Value = '  +4.24%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '7.128'
$ws.Range('E40').Value = '  +2.61%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1539'
$ws.Range('E41').Value = '  +1.88%  '
$ws.Range('B42').Value = 'Aptos'
$ws.Range('C42').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.411'
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('B43').Value = 'EnergySwap'
$ws.Range('C43').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.69'
$ws.Range('E43').Value = '  +2.11%  '
$ws.Range('B44').Value = 'Decentraland'
$ws.Range('C44').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.4844'
$ws.Range('E44').Value = '  +3.05%  '
$ws.Range('B45').Value = 'PaxDollar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.022'
$ws.Range('E45').Value = '  +1.37%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.677'
$ws.Range('E46').Value = '  +5.24%  '
$ws.Range('B47').Value = 'Quant'
$ws.Range('C47').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '103.48'
$ws.Range('E47').Value = '  +1.66%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '66.80'
$ws.Range('E48').Value = '  +2.56%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.06100'
$ws.Range('E49').Value = '  +1.06%  '
$ws.Range('B50').Value = 'EOS'
$ws.Range('C50').Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.9020'
$ws.Range('E50').Value = '  +1.83%  '
$ws.Range('B51').Value = 'Elrond'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '36.88'
$ws.Range('E51').Value = '  +1.69%  '
